# Adding the changes we made on may 9th
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Drop the two oldest gyroscope samples (old rows 2:3), shifting the rest up
$ws.Rows("2:3").Delete()

# Append the newly captured samples at the end of the series
$newSamples = @(
  @(-12.60414218902588, -15.56717586517334, 2.47767186164856),
  @(7.210841655731201, -12.07630348205566, 8.885769844055176),
  @(1.924999475479126, -10.15509986877441, 0.8653942346572876),
  @(5.971939086914063, 14.02054500579834, -4.391685962677002),
  @(-2.465487957000732, 0.488490343093872, 5.172791004180908),
  @(-6.064484119415283, 1.030177354812622, 4.91339921951294),
  @(-0.608199417591095, 12.21864986419678, -3.169827461242676),
  @(-0.2609232068061828, 3.39666223526001, -2.653707027435303),
  @(6.485929012298584, -1.756554484367371, -2.664892196655273),
  @(1.498893618583679, -2.116081237792969, -2.421479225158691),
  @(-0.8995492458343506, 1.140432238578796, 0.6262423396110535),
  @(-11.36417484283447, -11.00677871704102, -5.598630428314209)
)

$startRow = 20
for ($i = 0; $i -lt $newSamples.Count; $i++) {
  $row = $startRow + $i
  $ws.Cells.Item($row, 1).Value = $newSamples[$i][0]
  $ws.Cells.Item($row, 2).Value = $newSamples[$i][1]
  $ws.Cells.Item($row, 3).Value = $newSamples[$i][2]
}
